$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the 11.20.19 activity text: "documentation" -> "coumentation"
$ws.Cells.Item(9, 2).Value = "Added new features as a feature engineering task, algorithm testing code cleanup, initial writing of coumentation"

# Add a new log row 14 for 12.07.19 describing the final code run / doc writing
# Force column A to stay plain text (it would otherwise be auto-detected as a date)
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "12.07.19"
$ws.Cells.Item(14, 1).ClearFormats()

$ws.Cells.Item(14, 2).Value = "Finished notebook final code run, finished code part, refactored all notebook with new outline, started writing documentation"
$ws.Cells.Item(14, 2).WrapText = $true
$ws.Rows.Item(14).RowHeight = 29

# Update the active selection to B10 (matches the post-edit saved cursor position)
[void]$ws.Range("B10").Select()
